$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IBAN / account numbers with new test data
$ws.Range("B12").Value = "CH9689144715152235363"
$ws.Range("B15").Value = "DE03500105177284191473"
$ws.Range("B16").Value = "DE64500105177324151368"
$ws.Range("B17").Value = "IT37T0300203280682244199423"

# Update the active cell / selection to B17
$ws.Range("B17").Select()
